$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1050873333333333
$ws.Range("H2").Value = 0.315262
$ws.Range("I2").Value = 0.03031434174852429
$ws.Range("J2").Value = 0.03031434174852429
$ws.Range("M2").Value = 3.063353333333333
$ws.Range("N2").Value = 9.190059999999999
$ws.Range("O2").Value = 0.1884019917097105
$ws.Range("P2").Value = 0.1884019917097105
$ws.Range("Q2").Value = 0.3219196328577777
$ws.Range("R2").Value = 2.89727669572
$ws.Range("S2").Value = 0.005711282362790804
$ws.Range("T2").Value = 0.005711282362790804
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1050873333333333
$ws.Range("H3").Value = 0.315262
$ws.Range("I3").Value = 0.03031434174852429
$ws.Range("J3").Value = 0.03031434174852429
$ws.Range("O3").Value = 0.4156086771445645
$ws.Range("P3").Value = 0.4156086771445645
$ws.Range("Q3").Value = 0.7101442587986665
$ws.Range("R3").Value = 6.391298329187999
$ws.Range("S3").Value = 0.01259890347261242
$ws.Range("T3").Value = 0.01259890347261242
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.1050873333333333
$ws.Range("H4").Value = 0.315262
$ws.Range("I4").Value = 0.03031434174852429
$ws.Range("J4").Value = 0.03031434174852429
$ws.Range("M4").Value = 3.493414666666666
$ws.Range("N4").Value = 10.480244
$ws.Range("O4").Value = 0.214851572590793
$ws.Range("P4").Value = 0.214851572590793
$ws.Range("Q4").Value = 0.3671136315475555
$ws.Range("R4").Value = 3.304022683927999
$ws.Range("S4").Value = 0.006513083996725173
$ws.Range("T4").Value = 0.006513083996725173
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.1050873333333333
$ws.Range("H5").Value = 0.315262
$ws.Range("I5").Value = 0.03031434174852429
$ws.Range("J5").Value = 0.03031434174852429
$ws.Range("M5").Value = 2.945239333333333
$ws.Range("N5").Value = 8.835718
$ws.Range("O5").Value = 0.1811377585549322
$ws.Range("P5").Value = 0.1811377585549322
$ws.Range("Q5").Value = 0.3095073475684444
$ws.Range("R5").Value = 2.785566128116
$ws.Range("S5").Value = 0.005491071916395893
$ws.Range("T5").Value = 0.005491071916395893
$ws.Range("I6").Value = 0.6354599969768544
$ws.Range("J6").Value = 0.6354599969768545
$ws.Range("M6").Value = 3.063353333333333
$ws.Range("N6").Value = 9.190059999999999
$ws.Range("O6").Value = 0.1884019917097105
$ws.Range("P6").Value = 0.1884019917097105
$ws.Range("Q6").Value = 6.748193664226665
$ws.Range("R6").Value = 60.73374297803999
$ws.Range("S6").Value = 0.119721929082286
$ws.Range("T6").Value = 0.119721929082286
$ws.Range("I7").Value = 0.6354599969768544
$ws.Range("J7").Value = 0.6354599969768545
$ws.Range("O7").Value = 0.4156086771445645
$ws.Range("P7").Value = 0.4156086771445645
$ws.Range("S7").Value = 0.2641026887218394
$ws.Range("T7").Value = 0.2641026887218395
$ws.Range("I8").Value = 0.6354599969768544
$ws.Range("J8").Value = 0.6354599969768545
$ws.Range("M8").Value = 3.493414666666666
$ws.Range("N8").Value = 10.480244
$ws.Range("O8").Value = 0.214851572590793
$ws.Range("P8").Value = 0.214851572590793
$ws.Range("Q8").Value = 7.695566314077332
$ws.Range("R8").Value = 69.26009682669599
$ws.Range("S8").Value = 0.1365295796690177
$ws.Range("T8").Value = 0.1365295796690177
$ws.Range("I9").Value = 0.6354599969768544
$ws.Range("J9").Value = 0.6354599969768545
$ws.Range("M9").Value = 2.945239333333333
$ws.Range("N9").Value = 8.835718
$ws.Range("O9").Value = 0.1811377585549322
$ws.Range("P9").Value = 0.1811377585549322
$ws.Range("Q9").Value = 6.488002932134665
$ws.Range("R9").Value = 58.39202638921199
$ws.Range("S9").Value = 0.1151057995037114
$ws.Range("T9").Value = 0.1151057995037114
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.9666886666666668
$ws.Range("H10").Value = 2.900066
$ws.Range("I10").Value = 0.2788588279503266
$ws.Range("J10").Value = 0.2788588279503266
$ws.Range("M10").Value = 3.063353333333333
$ws.Range("N10").Value = 9.190059999999999
$ws.Range("O10").Value = 0.1884019917097105
$ws.Range("P10").Value = 0.1884019917097105
$ws.Range("Q10").Value = 2.961308949328889
$ws.Range("R10").Value = 26.65178054396
$ws.Range("S10").Value = 0.05253755859167702
$ws.Range("T10").Value = 0.05253755859167702
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.9666886666666668
$ws.Range("H11").Value = 2.900066
$ws.Range("I11").Value = 0.2788588279503266
$ws.Range("J11").Value = 0.2788588279503266
$ws.Range("O11").Value = 0.4156086771445645
$ws.Range("P11").Value = 0.4156086771445645
$ws.Range("Q11").Value = 6.532551401809333
$ws.Range("R11").Value = 58.792962616284
$ws.Range("S11").Value = 0.1158961485945189
$ws.Range("T11").Value = 0.1158961485945189
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.9666886666666668
$ws.Range("H12").Value = 2.900066
$ws.Range("I12").Value = 0.2788588279503266
$ws.Range("J12").Value = 0.2788588279503266
$ws.Range("M12").Value = 3.493414666666666
$ws.Range("N12").Value = 10.480244
$ws.Range("O12").Value = 0.214851572590793
$ws.Range("P12").Value = 0.214851572590793
$ws.Range("Q12").Value = 3.377044366233778
$ws.Range("R12").Value = 30.393399296104
$ws.Range("S12").Value = 0.05991325771595305
$ws.Range("T12").Value = 0.05991325771595305
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.9666886666666668
$ws.Range("H13").Value = 2.900066
$ws.Range("I13").Value = 0.2788588279503266
$ws.Range("J13").Value = 0.2788588279503266
$ws.Range("M13").Value = 2.945239333333333
$ws.Range("N13").Value = 8.835718
$ws.Range("O13").Value = 0.1811377585549322
$ws.Range("P13").Value = 0.1811377585549322
$ws.Range("Q13").Value = 2.847129484154222
$ws.Range("R13").Value = 25.624165357388
$ws.Range("S13").Value = 0.05051186304817763
$ws.Range("T13").Value = 0.05051186304817763
$ws.Range("G14").Value = 0.191934
$ws.Range("H14").Value = 0.575802
$ws.Range("I14").Value = 0.05536683332429468
$ws.Range("J14").Value = 0.05536683332429467
$ws.Range("M14").Value = 3.063353333333333
$ws.Range("N14").Value = 9.190059999999999
$ws.Range("O14").Value = 0.1884019917097105
$ws.Range("P14").Value = 0.1884019917097105
$ws.Range("Q14").Value = 0.5879616586800001
$ws.Range("R14").Value = 5.29165492812
$ws.Range("S14").Value = 0.01043122167295669
$ws.Range("T14").Value = 0.01043122167295669
$ws.Range("G15").Value = 0.191934
$ws.Range("H15").Value = 0.575802
$ws.Range("I15").Value = 0.05536683332429468
$ws.Range("J15").Value = 0.05536683332429467
$ws.Range("O15").Value = 0.4156086771445645
$ws.Range("P15").Value = 0.4156086771445645
$ws.Range("Q15").Value = 1.297024330572
$ws.Range("R15").Value = 11.673218975148
$ws.Range("S15").Value = 0.0230109363555937
$ws.Range("T15").Value = 0.0230109363555937
$ws.Range("G16").Value = 0.191934
$ws.Range("H16").Value = 0.575802
$ws.Range("I16").Value = 0.05536683332429468
$ws.Range("J16").Value = 0.05536683332429467
$ws.Range("M16").Value = 3.493414666666666
$ws.Range("N16").Value = 10.480244
$ws.Range("O16").Value = 0.214851572590793
$ws.Range("P16").Value = 0.214851572590793
$ws.Range("Q16").Value = 0.6705050506320001
$ws.Range("R16").Value = 6.034545455688
$ws.Range("S16").Value = 0.01189565120909703
$ws.Range("T16").Value = 0.01189565120909703
$ws.Range("G17").Value = 0.191934
$ws.Range("H17").Value = 0.575802
$ws.Range("I17").Value = 0.05536683332429468
$ws.Range("J17").Value = 0.05536683332429467
$ws.Range("M17").Value = 2.945239333333333
$ws.Range("N17").Value = 8.835718
$ws.Range("O17").Value = 0.1811377585549322
$ws.Range("P17").Value = 0.1811377585549322
$ws.Range("Q17").Value = 0.565291566204
$ws.Range("R17").Value = 5.087624095836
$ws.Range("S17").Value = 0.01002902408664726
$ws.Range("T17").Value = 0.01002902408664726
